# Fruta / hortaliza, semanal
# Insert two new weekly observation rows (new row 25 and 26) for
# "Vega Monumental Concepción - Arándano (blue)", pushing the existing
# rows 25-60 down to rows 27-62.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 25. Excel shifts all the
# rows below (25-60) down to (27-62), carrying their existing values,
# formatting and the date style on column D along with them.
$ws.Rows.Item(25).Insert()
$ws.Rows.Item(25).Insert()

# Fill in the data for the first new row (row 25).
$ws.Cells.Item(25, 1).Value  = 11
$ws.Cells.Item(25, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(25, 3).Value  = "Bíobío"
$ws.Cells.Item(25, 4).Value  = 44533
$ws.Cells.Item(25, 5).Value  = 8
$ws.Cells.Item(25, 6).Value  = "Fruta"
$ws.Cells.Item(25, 7).Value  = 100101
$ws.Cells.Item(25, 8).Value  = "Berries"
$ws.Cells.Item(25, 9).Value  = 100101001
$ws.Cells.Item(25, 10).Value = "Arándano (blue)"
$ws.Cells.Item(25, 11).Value = "Sin especificar"
$ws.Cells.Item(25, 12).Value = "Primera"
$ws.Cells.Item(25, 13).Value = 200
$ws.Cells.Item(25, 14).Value = 3800
$ws.Cells.Item(25, 15).Value = 4000
$ws.Cells.Item(25, 16).Value = 3900
$ws.Cells.Item(25, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(25, 18).Value = "Región de Ñuble"
$ws.Cells.Item(25, 19).Value = 1950
$ws.Cells.Item(25, 20).Value = 2

# Fill in the data for the second new row (row 26).
$ws.Cells.Item(26, 1).Value  = 11
$ws.Cells.Item(26, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(26, 3).Value  = "Bíobío"
$ws.Cells.Item(26, 4).Value  = 44533
$ws.Cells.Item(26, 5).Value  = 8
$ws.Cells.Item(26, 6).Value  = "Fruta"
$ws.Cells.Item(26, 7).Value  = 100101
$ws.Cells.Item(26, 8).Value  = "Berries"
$ws.Cells.Item(26, 9).Value  = 100101001
$ws.Cells.Item(26, 10).Value = "Arándano (blue)"
$ws.Cells.Item(26, 11).Value = "Sin especificar"
$ws.Cells.Item(26, 12).Value = "Segunda"
$ws.Cells.Item(26, 13).Value = 100
$ws.Cells.Item(26, 14).Value = 3500
$ws.Cells.Item(26, 15).Value = 3500
$ws.Cells.Item(26, 16).Value = 3500
$ws.Cells.Item(26, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(26, 18).Value = "Región de Ñuble"
$ws.Cells.Item(26, 19).Value = 1750
$ws.Cells.Item(26, 20).Value = 2
